# Apply the "Complex Combinations" redesign to the "Tests" sheet:
#  - Row 3 (N3): redesigned INDIRECT formula (dynamic reference via &)
#  - Row 4: new INDIRECT(...CHAR(66)...) formula plus companion O/P/Q cells
#  - Rows 1-2: new Aggregation (O), source-ref (P) and literal (Q/Z) helper cells
#  - Rows 5-7: additional P/Q helper cells
#  - Old row 20 marker cell is removed entirely (replaced by the richer grid)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# --- Row 1 -----------------------------------------------------------
# N1 is untouched: =INDEX(OFFSET(Data!A1, 0, 0, 3, 3), 2, 2)
$ws.Range("O1").Formula = "=SUM(INDEX(Data!A1:E6, 0, 2))"
$ws.Range("P1").Value = "Data!B2"
$ws.Range("Q1").Value = 25
$ws.Range("Z1").Value = "Test Value"

# --- Row 2 -------------------------------------------------------------
# N2 is untouched: =OFFSET(INDEX(Data!A1:E6, 2, 1), 1, 1)
$ws.Range("O2").Formula = "=AVERAGE(OFFSET(Data!B1, 1, 0, 5, 1))"
$ws.Range("P2").Value = "Data!C3"
$ws.Range("Q2").Value = "Bob"

# --- Row 3 --------------------------------------------------------------
# N3 formula redesigned to a simpler dynamic reference
$ws.Range("N3").Formula = "=INDIRECT(""Data!A"" & 2)"
$ws.Range("O3").Formula = "=COUNT(INDIRECT(""Data!B:B""))"
$ws.Range("P3").Value = "Data!A1:C3"
$ws.Range("Q3").Value = $true

# --- Row 4 (new) ---------------------------------------------------------
$ws.Range("N4").Formula = "=INDIRECT(""Data!"" & CHAR(66) & ""2"")"
$ws.Range("O4").Formula = "=MAX(INDEX(Data!A1:E6, 0, 4))"
$ws.Range("P4").Value = "InvalidSheet!A1"
$ws.Range("Q4").Value = "#REF!"

# --- Row 5 (new) -----------------------------------------------------------
# P5 is a blank placeholder cell (source data had no reference text here);
# an empty string assignment is a no-op in Excel (it clears/omits the cell),
# which is also the correct, real-Excel behaviour for a blank entry.
$ws.Range("P5").Value = ""
$ws.Range("Q5").Value = "#VALUE!"

# --- Row 6 (new) -----------------------------------------------------------
$ws.Range("P6").Value = "Data!A:A"

# --- Row 7 (new) -----------------------------------------------------------
$ws.Range("P7").Value = "Data!1:1"

# --- Remove the old row-20 marker cell entirely ---------------------------
$ws.Range("N20").ClearContents()
